$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Update the "Information" note text in row 4 (A4:J4 is a merged cell)
# ---------------------------------------------------------------------------
$ws.Range("A4").Value = "Information: excel file (duration of stay is the avg of the 12 months) and from the PDF for Spain (Scaled by rented accomodations)"

# ---------------------------------------------------------------------------
# Re-select row 5 (A5:J5) the way the author left the workbook
# ---------------------------------------------------------------------------
$ws.Range("A5:J5").Select()

# ---------------------------------------------------------------------------
# Country data rows (8-17): country name, bednights (D), avg length of stay (H)
# Arrivals (B) and nights-per-market (E) become formulas driven off D and H.
# ---------------------------------------------------------------------------
$rows = @(
    @{ Row = 8;  Country = "Germany";              D = 11176545; H = $null },
    @{ Row = 9;  Country = "Belgium";               D = 2538829;  H = $null },
    @{ Row = 10; Country = "France";                D = 11156671; H = $null },
    @{ Row = 11; Country = "Republic of Ireland";   D = 2177106;  H = 3.0792082197744302 },
    @{ Row = 12; Country = "Italy";                 D = 4542709;  H = 3.0792082197744302 },
    @{ Row = 13; Country = "Netherlands";           D = 3701944;  H = 3.0792082197744302 },
    @{ Row = 14; Country = "Portugal";               D = 2440746;  H = 3.0792082197744302 },
    @{ Row = 15; Country = "United Kingdom";        D = 18078076; H = 3.0792082197744302 },
    @{ Row = 16; Country = "Switzerland";           D = 1824839;  H = 3.0792082197744302 },
    @{ Row = 17; Country = "Russia";                D = 1311746;  H = 3.0792082197744302 }
)

foreach ($item in $rows) {
    $r = $item.Row

    # Country label in column A
    $ws.Range("A$r").Value = $item.Country

    # Column C is no longer used at all - remove the cell (value + formatting)
    $ws.Range("C$r").Clear()

    # Column D holds the raw bednights figure, formatted as a plain integer
    $ws.Range("D$r").Value = $item.D
    $ws.Range("D$r").NumberFormat = "#,##0"

    # Column B becomes a formula (scaled arrivals) with default/general formatting
    $ws.Range("B$r").Formula = "=D$r*(68314398/83701011)"
    $ws.Range("B$r").Style = "Normal"

    # Column E becomes a formula (nights) with an integer format and a top border
    $ws.Range("E$r").Formula = "=H$r*B$r"
    $ws.Range("E$r").NumberFormat = "#,##0"
    $ws.Range("E$r").Borders.Item(8).LineStyle = 1

    # Column H (average length of stay) only needs updating for the new rows
    if ($item.H -ne $null) {
        $ws.Range("H$r").Value = $item.H
    }
}
